# LMS-2523 Update BaSynthec Validation
# Update the "Header Format" row's example value on the openbis-metadata sheet
# from "MGP47" to "JJS-MGP47".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("openbis-metadata")
$ws.Range("B3").Value = "JJS-MGP47"
